$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 11).Value = $ws.Cells.Item(2, 10).Value2
$ws.Cells.Item(2, 10).Value = $ws.Cells.Item(2, 9).Value2
$ws.Cells.Item(2, 9).Value = $ws.Cells.Item(2, 8).Value2
$ws.Cells.Item(2, 8).Value = $ws.Cells.Item(2, 7).Value2
$ws.Cells.Item(2, 7).Value = $ws.Cells.Item(2, 6).Value2
$ws.Cells.Item(2, 6).Value = $ws.Cells.Item(2, 5).Value2
$ws.Cells.Item(2, 5).Value = $ws.Cells.Item(2, 4).Value2
$ws.Cells.Item(2, 4).Value = $ws.Cells.Item(2, 3).Value2
$ws.Cells.Item(2, 3).Value = $ws.Cells.Item(2, 2).Value2
$ws.Cells.Item(2, 2).Value = [double]"-3.965936795080616E-07"
$ws.Cells.Item(3, 11).Value = $ws.Cells.Item(3, 10).Value2
$ws.Cells.Item(3, 10).Value = $ws.Cells.Item(3, 9).Value2
$ws.Cells.Item(3, 9).Value = $ws.Cells.Item(3, 8).Value2
$ws.Cells.Item(3, 8).Value = $ws.Cells.Item(3, 7).Value2
$ws.Cells.Item(3, 7).Value = $ws.Cells.Item(3, 6).Value2
$ws.Cells.Item(3, 6).Value = $ws.Cells.Item(3, 5).Value2
$ws.Cells.Item(3, 5).Value = $ws.Cells.Item(3, 4).Value2
$ws.Cells.Item(3, 4).Value = $ws.Cells.Item(3, 3).Value2
$ws.Cells.Item(3, 3).Value = $ws.Cells.Item(3, 2).Value2
$ws.Cells.Item(3, 2).Value = [double]"-3.930720193778825E-10"
$ws.Cells.Item(4, 11).Value = $ws.Cells.Item(4, 10).Value2
$ws.Cells.Item(4, 10).Value = $ws.Cells.Item(4, 9).Value2
$ws.Cells.Item(4, 9).Value = $ws.Cells.Item(4, 8).Value2
$ws.Cells.Item(4, 8).Value = $ws.Cells.Item(4, 7).Value2
$ws.Cells.Item(4, 7).Value = $ws.Cells.Item(4, 6).Value2
$ws.Cells.Item(4, 6).Value = $ws.Cells.Item(4, 5).Value2
$ws.Cells.Item(4, 5).Value = $ws.Cells.Item(4, 4).Value2
$ws.Cells.Item(4, 4).Value = $ws.Cells.Item(4, 3).Value2
$ws.Cells.Item(4, 3).Value = $ws.Cells.Item(4, 2).Value2
$ws.Cells.Item(4, 2).Value = [double]"1.181302580199883E-07"
$ws.Cells.Item(5, 11).Value = $ws.Cells.Item(5, 10).Value2
$ws.Cells.Item(5, 10).Value = $ws.Cells.Item(5, 9).Value2
$ws.Cells.Item(5, 9).Value = $ws.Cells.Item(5, 8).Value2
$ws.Cells.Item(5, 8).Value = $ws.Cells.Item(5, 7).Value2
$ws.Cells.Item(5, 7).Value = $ws.Cells.Item(5, 6).Value2
$ws.Cells.Item(5, 6).Value = $ws.Cells.Item(5, 5).Value2
$ws.Cells.Item(5, 5).Value = $ws.Cells.Item(5, 4).Value2
$ws.Cells.Item(5, 4).Value = $ws.Cells.Item(5, 3).Value2
$ws.Cells.Item(5, 3).Value = $ws.Cells.Item(5, 2).Value2
$ws.Cells.Item(5, 2).Value = [double]"2.965444589886346E-07"
$ws.Cells.Item(6, 11).Value = $ws.Cells.Item(6, 10).Value2
$ws.Cells.Item(6, 10).Value = $ws.Cells.Item(6, 9).Value2
$ws.Cells.Item(6, 9).Value = $ws.Cells.Item(6, 8).Value2
$ws.Cells.Item(6, 8).Value = $ws.Cells.Item(6, 7).Value2
$ws.Cells.Item(6, 7).Value = $ws.Cells.Item(6, 6).Value2
$ws.Cells.Item(6, 6).Value = $ws.Cells.Item(6, 5).Value2
$ws.Cells.Item(6, 5).Value = $ws.Cells.Item(6, 4).Value2
$ws.Cells.Item(6, 4).Value = $ws.Cells.Item(6, 3).Value2
$ws.Cells.Item(6, 3).Value = $ws.Cells.Item(6, 2).Value2
$ws.Cells.Item(6, 2).Value = [double]"3.593882045849206E-07"
$ws.Cells.Item(7, 11).Value = $ws.Cells.Item(7, 10).Value2
$ws.Cells.Item(7, 10).Value = $ws.Cells.Item(7, 9).Value2
$ws.Cells.Item(7, 9).Value = $ws.Cells.Item(7, 8).Value2
$ws.Cells.Item(7, 8).Value = $ws.Cells.Item(7, 7).Value2
$ws.Cells.Item(7, 7).Value = $ws.Cells.Item(7, 6).Value2
$ws.Cells.Item(7, 6).Value = $ws.Cells.Item(7, 5).Value2
$ws.Cells.Item(7, 5).Value = $ws.Cells.Item(7, 4).Value2
$ws.Cells.Item(7, 4).Value = $ws.Cells.Item(7, 3).Value2
$ws.Cells.Item(7, 3).Value = $ws.Cells.Item(7, 2).Value2
$ws.Cells.Item(7, 2).Value = [double]"-1.035781544145298E-07"
$ws.Cells.Item(8, 11).Value = $ws.Cells.Item(8, 10).Value2
$ws.Cells.Item(8, 10).Value = $ws.Cells.Item(8, 9).Value2
$ws.Cells.Item(8, 9).Value = $ws.Cells.Item(8, 8).Value2
$ws.Cells.Item(8, 8).Value = $ws.Cells.Item(8, 7).Value2
$ws.Cells.Item(8, 7).Value = $ws.Cells.Item(8, 6).Value2
$ws.Cells.Item(8, 6).Value = $ws.Cells.Item(8, 5).Value2
$ws.Cells.Item(8, 5).Value = $ws.Cells.Item(8, 4).Value2
$ws.Cells.Item(8, 4).Value = $ws.Cells.Item(8, 3).Value2
$ws.Cells.Item(8, 3).Value = $ws.Cells.Item(8, 2).Value2
$ws.Cells.Item(8, 2).Value = [double]"-3.87512216759589E-10"
$ws.Cells.Item(9, 11).Value = $ws.Cells.Item(9, 10).Value2
$ws.Cells.Item(9, 10).Value = $ws.Cells.Item(9, 9).Value2
$ws.Cells.Item(9, 9).Value = $ws.Cells.Item(9, 8).Value2
$ws.Cells.Item(9, 8).Value = $ws.Cells.Item(9, 7).Value2
$ws.Cells.Item(9, 7).Value = $ws.Cells.Item(9, 6).Value2
$ws.Cells.Item(9, 6).Value = $ws.Cells.Item(9, 5).Value2
$ws.Cells.Item(9, 5).Value = $ws.Cells.Item(9, 4).Value2
$ws.Cells.Item(9, 4).Value = $ws.Cells.Item(9, 3).Value2
$ws.Cells.Item(9, 3).Value = $ws.Cells.Item(9, 2).Value2
$ws.Cells.Item(9, 2).Value = [double]"-1.07388789361007E-07"
$ws.Cells.Item(10, 11).Value = $ws.Cells.Item(10, 10).Value2
$ws.Cells.Item(10, 10).Value = $ws.Cells.Item(10, 9).Value2
$ws.Cells.Item(10, 9).Value = $ws.Cells.Item(10, 8).Value2
$ws.Cells.Item(10, 8).Value = $ws.Cells.Item(10, 7).Value2
$ws.Cells.Item(10, 7).Value = $ws.Cells.Item(10, 6).Value2
$ws.Cells.Item(10, 6).Value = $ws.Cells.Item(10, 5).Value2
$ws.Cells.Item(10, 5).Value = $ws.Cells.Item(10, 4).Value2
$ws.Cells.Item(10, 4).Value = $ws.Cells.Item(10, 3).Value2
$ws.Cells.Item(10, 3).Value = $ws.Cells.Item(10, 2).Value2
$ws.Cells.Item(10, 2).Value = [double]"-1.035472805832605E-07"
$ws.Cells.Item(11, 11).Value = $ws.Cells.Item(11, 10).Value2
$ws.Cells.Item(11, 10).Value = $ws.Cells.Item(11, 9).Value2
$ws.Cells.Item(11, 9).Value = $ws.Cells.Item(11, 8).Value2
$ws.Cells.Item(11, 8).Value = $ws.Cells.Item(11, 7).Value2
$ws.Cells.Item(11, 7).Value = $ws.Cells.Item(11, 6).Value2
$ws.Cells.Item(11, 6).Value = $ws.Cells.Item(11, 5).Value2
$ws.Cells.Item(11, 5).Value = $ws.Cells.Item(11, 4).Value2
$ws.Cells.Item(11, 4).Value = $ws.Cells.Item(11, 3).Value2
$ws.Cells.Item(11, 3).Value = $ws.Cells.Item(11, 2).Value2
$ws.Cells.Item(11, 2).Value = [double]"6.303355340908645E-06"
$ws.Cells.Item(12, 10).Value = $ws.Cells.Item(12, 9).Value2
$ws.Cells.Item(12, 9).Value = $ws.Cells.Item(12, 8).Value2
$ws.Cells.Item(12, 8).Value = $ws.Cells.Item(12, 7).Value2
$ws.Cells.Item(12, 7).Value = $ws.Cells.Item(12, 6).Value2
$ws.Cells.Item(12, 6).Value = $ws.Cells.Item(12, 5).Value2
$ws.Cells.Item(12, 5).Value = $ws.Cells.Item(12, 4).Value2
$ws.Cells.Item(12, 4).Value = $ws.Cells.Item(12, 3).Value2
$ws.Cells.Item(12, 3).Value = $ws.Cells.Item(12, 2).Value2
$ws.Cells.Item(12, 2).Value = [double]"-2.375649628613696E-07"
$ws.Cells.Item(13, 9).Value = $ws.Cells.Item(13, 8).Value2
$ws.Cells.Item(13, 8).Value = $ws.Cells.Item(13, 7).Value2
$ws.Cells.Item(13, 7).Value = $ws.Cells.Item(13, 6).Value2
$ws.Cells.Item(13, 6).Value = $ws.Cells.Item(13, 5).Value2
$ws.Cells.Item(13, 5).Value = $ws.Cells.Item(13, 4).Value2
$ws.Cells.Item(13, 4).Value = $ws.Cells.Item(13, 3).Value2
$ws.Cells.Item(13, 3).Value = $ws.Cells.Item(13, 2).Value2
$ws.Cells.Item(13, 2).Value = [double]"3.720025918141356E-07"
$ws.Cells.Item(14, 8).Value = $ws.Cells.Item(14, 7).Value2
$ws.Cells.Item(14, 7).Value = $ws.Cells.Item(14, 6).Value2
$ws.Cells.Item(14, 6).Value = $ws.Cells.Item(14, 5).Value2
$ws.Cells.Item(14, 5).Value = $ws.Cells.Item(14, 4).Value2
$ws.Cells.Item(14, 4).Value = $ws.Cells.Item(14, 3).Value2
$ws.Cells.Item(14, 3).Value = $ws.Cells.Item(14, 2).Value2
$ws.Cells.Item(14, 2).Value = [double]"3.829984367986761E-07"
$ws.Cells.Item(15, 7).Value = $ws.Cells.Item(15, 6).Value2
$ws.Cells.Item(15, 6).Value = $ws.Cells.Item(15, 5).Value2
$ws.Cells.Item(15, 5).Value = $ws.Cells.Item(15, 4).Value2
$ws.Cells.Item(15, 4).Value = $ws.Cells.Item(15, 3).Value2
$ws.Cells.Item(15, 3).Value = $ws.Cells.Item(15, 2).Value2
$ws.Cells.Item(15, 2).Value = [double]"-3.160475492397508E-06"
$ws.Cells.Item(16, 6).Value = $ws.Cells.Item(16, 5).Value2
$ws.Cells.Item(16, 5).Value = $ws.Cells.Item(16, 4).Value2
$ws.Cells.Item(16, 4).Value = $ws.Cells.Item(16, 3).Value2
$ws.Cells.Item(16, 3).Value = $ws.Cells.Item(16, 2).Value2
$ws.Cells.Item(16, 2).Value = [double]"-4.101096154340844E-08"
$ws.Cells.Item(17, 5).Value = $ws.Cells.Item(17, 4).Value2
$ws.Cells.Item(17, 4).Value = $ws.Cells.Item(17, 3).Value2
$ws.Cells.Item(17, 3).Value = $ws.Cells.Item(17, 2).Value2
$ws.Cells.Item(17, 2).Value = [double]"-1.831659499074156E-07"
$ws.Cells.Item(18, 4).Value = $ws.Cells.Item(18, 3).Value2
$ws.Cells.Item(18, 3).Value = $ws.Cells.Item(18, 2).Value2
$ws.Cells.Item(18, 2).Value = [double]"2.770877186031306E-07"
$ws.Cells.Item(19, 3).Value = $ws.Cells.Item(19, 2).Value2
$ws.Cells.Item(19, 2).Value = [double]"2.29775004800814E-07"
$ws.Cells.Item(20, 2).Value = [double]"-1.554241066958895E-07"
